# Auto-generated edit script: update '想去人数' (column F) values
# across all four worksheets per the source diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 6686  # F2: 6672 -> 6686
$ws.Cells.Item(3, 6).Value = 795  # F3: 791 -> 795
$ws.Cells.Item(5, 6).Value = 128  # F5: 127 -> 128
$ws.Cells.Item(6, 6).Value = 697  # F6: 693 -> 697
$ws.Cells.Item(7, 6).Value = 219  # F7: 217 -> 219
$ws.Cells.Item(8, 6).Value = 10  # F8: 8 -> 10
$ws.Cells.Item(9, 6).Value = 1070  # F9: 1068 -> 1070
$ws.Cells.Item(10, 6).Value = 838  # F10: 833 -> 838
$ws.Cells.Item(12, 6).Value = 1317  # F12: 1313 -> 1317
$ws.Cells.Item(14, 6).Value = 109  # F14: 107 -> 109
$ws.Cells.Item(16, 6).Value = 539  # F16: 537 -> 539
$ws.Cells.Item(19, 6).Value = 1053  # F19: 1052 -> 1053
$ws.Cells.Item(20, 6).Value = 1472  # F20: 1470 -> 1472
$ws.Cells.Item(21, 6).Value = 713  # F21: 710 -> 713
$ws.Cells.Item(22, 6).Value = 119  # F22: 30 -> 119
$ws.Cells.Item(23, 6).Value = 445  # F23: 442 -> 445
$ws.Cells.Item(27, 6).Value = 1121  # F27: 1120 -> 1121
$ws.Cells.Item(28, 6).Value = 242  # F28: 241 -> 242
$ws.Cells.Item(29, 6).Value = 2339  # F29: 2338 -> 2339
$ws.Cells.Item(30, 6).Value = 262  # F30: 260 -> 262
$ws.Cells.Item(31, 6).Value = 1189  # F31: 1186 -> 1189
$ws.Cells.Item(34, 6).Value = 3808  # F34: 3800 -> 3808
$ws.Cells.Item(36, 6).Value = 693  # F36: 692 -> 693

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 1032  # F9: 1031 -> 1032
$ws.Cells.Item(12, 6).Value = 640  # F12: 639 -> 640
$ws.Cells.Item(17, 6).Value = 393  # F17: 391 -> 393
$ws.Cells.Item(18, 6).Value = 330  # F18: 329 -> 330
$ws.Cells.Item(19, 6).Value = 4120  # F19: 4118 -> 4120
$ws.Cells.Item(24, 6).Value = 1  # F24: 0 -> 1
$ws.Cells.Item(25, 6).Value = 229  # F25: 228 -> 229
$ws.Cells.Item(26, 6).Value = 244  # F26: 242 -> 244

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 1240  # F4: 1238 -> 1240
$ws.Cells.Item(5, 6).Value = 1623  # F5: 1621 -> 1623
$ws.Cells.Item(8, 6).Value = 937  # F8: 936 -> 937

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1240  # F2: 1238 -> 1240
$ws.Cells.Item(3, 6).Value = 1623  # F3: 1621 -> 1623
$ws.Cells.Item(6, 6).Value = 937  # F6: 936 -> 937
$ws.Cells.Item(9, 6).Value = 6686  # F9: 6672 -> 6686
$ws.Cells.Item(12, 6).Value = 795  # F12: 791 -> 795
$ws.Cells.Item(14, 6).Value = 129  # F14: 127 -> 129
$ws.Cells.Item(15, 6).Value = 697  # F15: 693 -> 697
$ws.Cells.Item(16, 6).Value = 1070  # F16: 1068 -> 1070
$ws.Cells.Item(17, 6).Value = 838  # F17: 833 -> 838
$ws.Cells.Item(22, 6).Value = 1317  # F22: 1313 -> 1317
$ws.Cells.Item(25, 6).Value = 539  # F25: 537 -> 539
$ws.Cells.Item(26, 6).Value = 330  # F26: 329 -> 330
$ws.Cells.Item(28, 6).Value = 1472  # F28: 1470 -> 1472
$ws.Cells.Item(29, 6).Value = 713  # F29: 710 -> 713
$ws.Cells.Item(30, 6).Value = 445  # F30: 442 -> 445
$ws.Cells.Item(32, 6).Value = 229  # F32: 228 -> 229
$ws.Cells.Item(33, 6).Value = 244  # F33: 242 -> 244
$ws.Cells.Item(34, 6).Value = 1121  # F34: 1120 -> 1121
$ws.Cells.Item(35, 6).Value = 242  # F35: 241 -> 242
$ws.Cells.Item(38, 6).Value = 2339  # F38: 2338 -> 2339
$ws.Cells.Item(45, 6).Value = 1189  # F45: 1186 -> 1189
$ws.Cells.Item(48, 6).Value = 3808  # F48: 3800 -> 3808
$ws.Cells.Item(50, 6).Value = 693  # F50: 692 -> 693

